# Template.xlsx update:
#  - remove the stray defined name "_xlnm._FilterDatabase_0_0"
#  - change the frozen-pane / view template from freezing columns A:B
#    (split at C2) to freezing only column A (split at B2), and move the
#    active selection in the bottom-right pane from A6 to A10
#  - update the per-store "Cross-Category equipment snacks" (column C)
#    display counts: clear some stores that no longer get that display
#    and set/raise the count for stores that now do

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# --- 1. Drop the duplicate/stray defined name ---------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase_0_0") {
        $n.Delete()
    }
}

# --- 2. Re-point the frozen panes / selection (new display template) ----
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$ws.Range("A10").Select()

# --- 3. Update column C ("Cross-Category equipment snacks") counts ------

# Stores that lose the cross-category display (count cleared)
$clearRows = @(5,6,16,21,24,25,26,28,29,30,31,32,35,38,39,40,44,49)
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 3).ClearContents()
}

# Stores that gain / change the cross-category display count
$setRows = [ordered]@{
    55 = 1; 62 = 1; 64 = 1; 66 = 1; 67 = 1; 68 = 1; 70 = 1; 71 = 1; 72 = 1;
    73 = 2; 74 = 1; 75 = 2; 78 = 1; 81 = 1; 82 = 1; 84 = 1; 89 = 2; 97 = 1
}
foreach ($r in $setRows.Keys) {
    $ws.Cells.Item($r, 3).Value = $setRows[$r]
}
